$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old "Sign-In" cell content (A1) entirely
$ws.Range("A1").ClearContents()

# Update A2 to the new product string
$ws.Range("A2").Value = "Motorola edge+ 5G UW"

# Adjust column A width (the engine stores width as ColumnWidth + 5/6,
# matching Excel's internal MDW-based padding for the default font;
# 21 + 1/6 round-trips to a stored width of exactly 22)
$ws.Columns.Item(1).ColumnWidth = 21.16666666666667

# Select the whole column A
$ws.Range("A1:A1048576").Select()
